# Ajuste Retanqueo (Carteras) - Digicredito
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RetanqueoDigiCreditoCCS")
$ws.Activate()

# Update the sample data row (row 2) with the new Digicredito test values
$ws.Range("B2").Value = '"10092369"'
$ws.Range("C2").Value = '"134950"'
$ws.Range("F2").Value = '"50"'
$ws.Range("G2").Value = '"50"'
$ws.Range("H2").Value = '"8700000"'
$ws.Range("I2").Value = '"200000"'
$ws.Range("J2").Value = '"300000"'
$ws.Range("L2").Value = '"450000"'
$ws.Range("M2").Value = '"OSCAR"'
$ws.Range("Q2").Value = '"200000"'
$ws.Range("R2").Value = '"250000"'
$ws.Range("S2").Value = '"3228483385"'
$ws.Range("V2").Value = '"25/11/2021"'
$ws.Range("X2").Value = '"86308"'

# Move the active selection to match the saved cursor position
$ws.Range("R8").Select()
